$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039822463303315
$ws.Range("D2").Value = 1.042121712432954
$ws.Range("E2").Value = 1.047808065905107
$ws.Range("F2").Value = 1.057498777878844
$ws.Range("I2").Value = 1.042918224583322
$ws.Range("J2").Value = 1.044912398873432
$ws.Range("K2").Value = 1.044899293951857
$ws.Range("L2").Value = 1.050569669777073
$ws.Range("M2").Value = 1.060233581550296
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040682827700167
$ws.Range("D3").Value = 1.042759306656269
$ws.Range("E3").Value = 1.048585319699211
$ws.Range("F3").Value = 1.058386113757658
$ws.Range("I3").Value = 1.043150409370335
$ws.Range("J3").Value = 1.045418573047117
$ws.Range("K3").Value = 1.045348109437201
$ws.Range("L3").Value = 1.051158932055791
$ws.Range("M3").Value = 1.060934578931261
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.041240068048998
$ws.Range("D4").Value = 1.043172236401727
$ws.Range("E4").Value = 1.049089122854325
$ws.Range("F4").Value = 1.058961317616562
$ws.Range("I4").Value = 1.043299591956946
$ws.Range("J4").Value = 1.04574596631244
$ws.Range("K4").Value = 1.045638197542666
$ws.Range("L4").Value = 1.051540431010172
$ws.Range("M4").Value = 1.061388569452126
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04147445656412
$ws.Range("D5").Value = 1.04334591746405
$ws.Range("E5").Value = 1.049301127768389
$ws.Range("F5").Value = 1.059203380062493
$ws.Range("I5").Value = 1.043362054614016
$ws.Range("J5").Value = 1.045883568844992
$ws.Range("K5").Value = 1.045760071459706
$ws.Range("L5").Value = 1.05170086132877
$ws.Range("M5").Value = 1.061579520888876
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.041513818694332
$ws.Range("D6").Value = 1.043375084238061
$ws.Range("E6").Value = 1.049336736359991
$ws.Range("F6").Value = 1.059244037828963
$ws.Range("I6").Value = 1.043372527469259
$ws.Range("J6").Value = 1.045906670910369
$ws.Range("K6").Value = 1.045780529952234
$ws.Range("L6").Value = 1.051727801069559
$ws.Range("M6").Value = 1.061611587932272
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04124319947181
$ws.Range("D7").Value = 1.043174556803119
$ws.Range("E7").Value = 1.049091954866828
$ws.Range("F7").Value = 1.058964551099543
$ws.Range("I7").Value = 1.04330042758387
$ws.Range("J7").Value = 1.045747805097744
$ws.Range("K7").Value = 1.045639826340302
$ws.Range("L7").Value = 1.051542574500424
$ws.Range("M7").Value = 1.061391120587847
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040113117443816
$ws.Range("D8").Value = 1.042337114092236
$ws.Range("E8").Value = 1.048070561704962
$ws.Range("F8").Value = 1.05779844113905
$ws.Range("I8").Value = 1.042996910807385
$ws.Range("J8").Value = 1.045083490085615
$ws.Range("K8").Value = 1.045051039932464
$ws.Range("L8").Value = 1.050768770071915
$ws.Range("M8").Value = 1.060470403490928
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038125871715059
$ws.Range("D9").Value = 1.040864291502418
$ws.Range("E9").Value = 1.046277460487591
$ws.Range("F9").Value = 1.055751637725469
$ws.Range("I9").Value = 1.042454018258114
$ws.Range("J9").Value = 1.04391189969904
$ws.Range("K9").Value = 1.044011082560252
$ws.Range("L9").Value = 1.049406872369339
$ws.Range("M9").Value = 1.058851101320852
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036803889301303
$ws.Range("D10").Value = 1.039884430488776
$ws.Range("E10").Value = 1.045086683406827
$ws.Range("F10").Value = 1.054392609441817
$ws.Range("I10").Value = 1.042086718776604
$ws.Range("J10").Value = 1.043130250922796
$ws.Range("K10").Value = 1.043316207619356
$ws.Range("L10").Value = 1.048500127417078
$ws.Range("M10").Value = 1.057773755320334
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036232150528977
$ws.Range("D11").Value = 1.039460639974644
$ws.Range("E11").Value = 1.044572182147071
$ws.Range("F11").Value = 1.053805464082878
$ws.Range("I11").Value = 1.04192640975119
$ws.Range("J11").Value = 1.042791663911964
$ws.Range("K11").Value = 1.043014961594752
$ws.Range("L11").Value = 1.048107795411837
$ws.Range("M11").Value = 1.057307791052855
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036019886192581
$ws.Range("D12").Value = 1.039303301472835
$ws.Range("E12").Value = 1.044381242448094
$ws.Range("F12").Value = 1.053587572721888
$ws.Range("I12").Value = 1.04186667440709
$ws.Range("J12").Value = 1.042665879606897
$ws.Range("K12").Value = 1.042903012421056
$ws.Range("L12").Value = 1.047962111456734
$ws.Range("M12").Value = 1.057134792901744
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036065412842439
$ws.Range("D13").Value = 1.039337047642014
$ws.Range("E13").Value = 1.044422191987024
$ws.Range("F13").Value = 1.053634302053059
$ws.Range("I13").Value = 1.041879496401308
$ws.Range("J13").Value = 1.042692861563095
$ws.Range("K13").Value = 1.042927028277552
$ws.Range("L13").Value = 1.047993359073255
$ws.Range("M13").Value = 1.057171897869322
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036214602544176
$ws.Range("D14").Value = 1.039447632761562
$ws.Range("E14").Value = 1.044556395554789
$ws.Range("F14").Value = 1.053787449012931
$ws.Range("I14").Value = 1.041921475869031
$ws.Range("J14").Value = 1.042781266902561
$ws.Range("K14").Value = 1.043005708915187
$ws.Range("L14").Value = 1.048095752192108
$ws.Range("M14").Value = 1.057293489295596
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036306537188105
$ws.Range("D15").Value = 1.039515778046232
$ws.Range("E15").Value = 1.044639105245606
$ws.Range("F15").Value = 1.053881834541162
$ws.Range("I15").Value = 1.041947315725378
$ws.Range("J15").Value = 1.042835733999072
$ws.Range("K15").Value = 1.043054179667023
$ws.Range("L15").Value = 1.048158846058094
$ws.Range("M15").Value = 1.057368416662018
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036841848349829
$ws.Range("D16").Value = 1.039912566658297
$ws.Range("E16").Value = 1.045120852761581
$ws.Range("F16").Value = 1.054431604372295
$ws.Range("I16").Value = 1.042097331351604
$ws.Range("J16").Value = 1.0431527192352
$ws.Range("K16").Value = 1.043336192830019
$ws.Range("L16").Value = 1.04852617152264
$ws.Range("M16").Value = 1.057804691209919
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037177820233702
$ws.Range("D17").Value = 1.040161595590589
$ws.Range("E17").Value = 1.04542333948814
$ws.Range("F17").Value = 1.05477681591607
$ws.Range("I17").Value = 1.042191093799891
$ws.Range("J17").Value = 1.043351522299749
$ws.Range("K17").Value = 1.043512996534572
$ws.Range("L17").Value = 1.048756664905436
$ws.Range("M17").Value = 1.058078498676515
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037373853264066
$ws.Range("D18").Value = 1.040306897675695
$ws.Range("E18").Value = 1.045599882156257
$ws.Range("F18").Value = 1.05497829948313
$ws.Range("I18").Value = 1.042245661655281
$ws.Range("J18").Value = 1.043467468339224
$ws.Range("K18").Value = 1.043616088272854
$ws.Range("L18").Value = 1.0488911361233
$ws.Range("M18").Value = 1.058238257325912
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03744070666445
$ws.Range("D19").Value = 1.040356449975116
$ws.Range("E19").Value = 1.045660096786794
$ws.Range("F19").Value = 1.055047021745821
$ws.Range("I19").Value = 1.042264247134528
$ws.Range("J19").Value = 1.043507000798235
$ws.Range("K19").Value = 1.043651233940576
$ws.Range("L19").Value = 1.048936992096715
$ws.Range("M19").Value = 1.058292739565779
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037141766760198
$ws.Range("D20").Value = 1.040134872194685
$ws.Range("E20").Value = 1.045390874416166
$ws.Range("F20").Value = 1.054739764781927
$ws.Range("I20").Value = 1.042181046610223
$ws.Range("J20").Value = 1.043330193888797
$ws.Range("K20").Value = 1.043494030768657
$ws.Range("L20").Value = 1.048731932206922
$ws.Range("M20").Value = 1.058049116400302
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036170666987761
$ws.Range("D21").Value = 1.039415066084122
$ws.Range("E21").Value = 1.044516871270373
$ws.Range("F21").Value = 1.05374234549903
$ws.Range("I21").Value = 1.041909119185325
$ws.Range("J21").Value = 1.042755234221105
$ws.Range("K21").Value = 1.042982540881637
$ws.Range("L21").Value = 1.048065598700355
$ws.Range("M21").Value = 1.057257681382699
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035560705676821
$ws.Range("D22").Value = 1.038962936907341
$ws.Range("E22").Value = 1.043968329318528
$ws.Range("F22").Value = 1.053116389736099
$ws.Range("I22").Value = 1.041737051988008
$ws.Range("J22").Value = 1.042393630634233
$ws.Range("K22").Value = 1.042660640908983
$ws.Range("L22").Value = 1.047646913000172
$ws.Range("M22").Value = 1.056760548377171
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03588399961259
$ws.Range("D23").Value = 1.039202576729481
$ws.Range("E23").Value = 1.044259028449818
$ws.Range("F23").Value = 1.053448110039895
$ws.Range("I23").Value = 1.041828371665051
$ws.Range("J23").Value = 1.042585332956637
$ws.Range("K23").Value = 1.042831314781901
$ws.Range("L23").Value = 1.047868840671514
$ws.Range("M23").Value = 1.05702404258024
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037158057568233
$ws.Range("D24").Value = 1.040146947199755
$ws.Range("E24").Value = 1.045405543653492
$ws.Range("F24").Value = 1.054756506200822
$ws.Range("I24").Value = 1.042185586880402
$ws.Range("J24").Value = 1.043339831318166
$ws.Range("K24").Value = 1.043502600689789
$ws.Range("L24").Value = 1.04874310776141
$ws.Range("M24").Value = 1.058062392827585
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038639126929203
$ws.Range("D25").Value = 1.041244702509481
$ws.Range("E25").Value = 1.046740212697546
$ws.Range("F25").Value = 1.056279822639949
$ws.Range("I25").Value = 1.042595319072201
$ws.Range("J25").Value = 1.044214892902563
$ws.Range("K25").Value = 1.044280218381713
$ws.Range("L25").Value = 1.049758752320429
$ws.Range("M25").Value = 1.059269350678369
